$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'63.126.29"
$ws.Range("E2").Value = "  -1.67%  "

# Row 3
$ws.Range("D3").Value = "'3.098.92"
$ws.Range("E3").Value = "  -1.30%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "'591.84"
$ws.Range("E5").Value = "  -3.41%  "

# Row 6
$ws.Range("D6").Value = "'136.80"
$ws.Range("E6").Value = "  -4.02%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").Value = "'3.093.57"
$ws.Range("E8").Value = "  -1.43%  "

# Row 9
$ws.Range("D9").Value = "'0.519"
$ws.Range("E9").Value = "  -0.33%  "

# Row 10
$ws.Range("D10").Value = "'0.146"
$ws.Range("E10").Value = "  -1.89%  "

# Row 11
$ws.Range("E11").Value = "  -0.37%  "

# Row 12
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  -2.05%  "

# Row 13
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("E13").Value = "  -1.93%  "

# Row 14
$ws.Range("D14").Value = "'34.03"
$ws.Range("E14").Value = "  -3.27%  "

# Row 15
$ws.Range("D15").Value = "'3.624.11"
$ws.Range("E15").Value = "  -1.13%  "

# Row 16
$ws.Range("E16").Value = "  +2.08%  "

# Row 17
$ws.Range("D17").Value = "'63.321.54"
$ws.Range("E17").Value = "  -1.43%  "

# Row 18
$ws.Range("D18").Value = "'3.125.86"
$ws.Range("E18").Value = "  -0.61%  "

# Row 19
$ws.Range("D19").Value = "'6.68"
$ws.Range("E19").Value = "  -1.54%  "

# Row 20
$ws.Range("D20").Value = "'479.49"
$ws.Range("E20").Value = "  +1.32%  "

# Row 21
$ws.Range("D21").Value = "'14.31"
$ws.Range("E21").Value = "  -1.42%  "

# Row 22
$ws.Range("D22").Value = "'0.694"
$ws.Range("E22").Value = "  -3.52%  "

# Row 23
$ws.Range("D23").Value = "'7.60"
$ws.Range("E23").Value = "  -2.96%  "

# Row 24
$ws.Range("D24").Value = "'87.10"
$ws.Range("E24").Value = "  +3.37%  "

# Row 25
$ws.Range("D25").Value = "'12.95"
$ws.Range("E25").Value = "  -4.47%  "

# Row 26
$ws.Range("E26").Value = "  +0.04%  "

# Row 27
$ws.Range("D27").Value = "'2.72"
$ws.Range("E27").Value = "  -2.63%  "

# Row 28
$ws.Range("D28").Value = "'8.05"
$ws.Range("E28").Value = "  -4.61%  "

# Row 29
$ws.Range("D29").Value = "'6.93"
$ws.Range("E29").Value = "  -2.64%  "

# Row 30
$ws.Range("D30").Value = "'2.03"
$ws.Range("E30").Value = "  -2.45%  "

# Row 31
$ws.Range("D31").Value = "'27.15"
$ws.Range("E31").Value = "  +3.69%  "

# Row 32
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.20%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.109"
$ws.Range("E33").Value = "  -10.81%  "

# Row 34
$ws.Range("D34").Value = "'2.54"
$ws.Range("E34").Value = "  -2.53%  "

# Row 35
$ws.Range("D35").Value = "'1.08"
$ws.Range("E35").Value = "  -2.88%  "

# Row 36
$ws.Range("D36").Value = "'5.91"
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("D37").Value = "'52.38"
$ws.Range("E37").Value = "  -0.66%  "

# Row 38
$ws.Range("D38").Value = "'0.0₃0718"
$ws.Range("E38").Value = "  -5.98%  "

# Row 39
$ws.Range("D39").Value = "'0.0388"
$ws.Range("E39").Value = "  -0.84%  "

# Row 40
$ws.Range("D40").Value = "'418.90"
$ws.Range("E40").Value = "  -7.30%  "

# Row 41
$ws.Range("D41").Value = "'2.78"
$ws.Range("E41").Value = "  -9.37%  "

# Row 42
$ws.Range("D42").Value = "'0.118"
$ws.Range("E42").Value = "  +0.15%  "

# Row 43
$ws.Range("D43").Value = "'8.24"
$ws.Range("E43").Value = "  -0.10%  "

# Row 44
$ws.Range("D44").Value = "'2.851.34"
$ws.Range("E44").Value = "  +1.06%  "

# Row 45
$ws.Range("D45").Value = "'0.256"
$ws.Range("E45").Value = "  -2.47%  "

# Row 46
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  -0.13%  "

# Row 47
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.13"
$ws.Range("E47").Value = "  -5.70%  "

# Row 48
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "'2.34"
$ws.Range("E48").Value = "  -2.91%  "

# Row 49
$ws.Range("D49").Value = "'0.113"
$ws.Range("E49").Value = "  +0.11%  "

# Row 50
$ws.Range("D50").Value = "'25.35"
$ws.Range("E50").Value = "  -3.32%  "

# Row 51
$ws.Range("D51").Value = "'120.07"
$ws.Range("E51").Value = "  +0.05%  "
